# 自动更新Excel文件 - 2026-02-12 23:21:11
# For every data row, the "剩余" (remaining days, column E) counter ticks
# down by 1 day. When a row's remaining count has reached 1 (i.e. it is
# about to expire), it is restocked/renewed: the remaining count resets
# to 10 and the "开始时间" (start date, column F) is pushed forward by
# 10 days. Rows whose start date is not a valid yyyyMMdd date are left
# untouched (stale/broken data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 99
}

for ($r = 2; $r -le $lastRow; $r++) {
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($eVal -eq $null -or $fVal -eq $null) {
        continue
    }

    $fStr = [string]([int]$fVal)

    $parsedDate = $null
    if ($fStr.Length -eq 8) {
        try {
            $parsedDate = [datetime]::ParseExact($fStr, "yyyyMMdd", $null)
        } catch {
            $parsedDate = $null
        }
    }

    if ($parsedDate -eq $null) {
        # invalid / unparsable start date -> leave row unchanged
        continue
    }

    if ([int]$eVal -eq 1) {
        $eCell.Value2 = 10
        $newDate = $parsedDate.AddDays(10)
        $fCell.Value2 = [int]$newDate.ToString("yyyyMMdd")
    } else {
        $eCell.Value2 = [int]$eVal - 1
    }
}
